$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("C8").Value = [double]"1.848651158265418e-08"
$ws.Range("D8").Value = 0
$ws.Range("C9").Value = [double]"8.629401430334518e-06"
$ws.Range("D9").Value = 0
$ws.Range("C10").Value = [double]"3.57861057156248e-05"
$ws.Range("D10").Value = 0
$ws.Range("C11").Value = [double]"8.158993579336011e-05"
$ws.Range("D11").Value = 0
$ws.Range("C12").Value = 0.0001452887697347207
$ws.Range("D12").Value = 0
$ws.Range("C13").Value = 0.0002264599827342399
$ws.Range("D13").Value = 0
$ws.Range("C14").Value = 0.0003247913490243328
$ws.Range("D14").Value = 0
$ws.Range("C15").Value = 0.0004408157094592831
$ws.Range("D15").Value = 0
$ws.Range("C16").Value = 0.0005735910967008442
$ws.Range("D16").Value = [double]"5.581914313972155e-06"
$ws.Range("C17").Value = 0.000724910568929212
$ws.Range("D17").Value = [double]"2.395874574452934e-05"
$ws.Range("C18").Value = 0.0008940871743036786
$ws.Range("D18").Value = [double]"5.395803128339271e-05"
$ws.Range("C19").Value = 0.001082678300202107
$ws.Range("D19").Value = [double]"9.411622673176995e-05"
$ws.Range("C20").Value = 0.001289899307701796
$ws.Range("D20").Value = 0.0001449428268263245
$ws.Range("C21").Value = 0.001515372281024599
$ws.Range("D21").Value = 0.0002058285475018214
$ws.Range("C22").Value = 0.001758703238329739
$ws.Range("D22").Value = 0.0002769013941261135
$ws.Range("C23").Value = 0.002021123038553598
$ws.Range("D23").Value = 0.0003581832471944157
$ws.Range("C24").Value = 0.002301638373605559
$ws.Range("D24").Value = 0.0004501173684320405
$ws.Range("C25").Value = 0.002601505142314592
$ws.Range("D25").Value = 0.0005524386362640647
$ws.Range("C26").Value = 0.002922015500598368
$ws.Range("D26").Value = 0.0006658570052580817
$ws.Range("C27").Value = 0.003262769751446372
$ws.Range("D27").Value = 0.0007898883685609782
$ws.Range("C28").Value = 0.003623793476160117
$ws.Range("D28").Value = 0.0009258213359183565
$ws.Range("C29").Value = 0.004006090849688459
$ws.Range("D29").Value = 0.001073456657501028
$ws.Range("C30").Value = 0.004409755602286333
$ws.Range("D30").Value = 0.001234174215762358
$ws.Range("C31").Value = 0.004835587592553756
$ws.Range("D31").Value = 0.001406994624885992
$ws.Range("C32").Value = 0.005283813535389457
$ws.Range("D32").Value = 0.001591301572885864
$ws.Range("C33").Value = 0.00575586740423954
$ws.Range("D33").Value = 0.001787480436637567
$ws.Range("C34").Value = 0.006252026569685305
$ws.Range("D34").Value = 0.001996563323604301
$ws.Range("C35").Value = 0.006771480875587036
$ws.Range("D35").Value = 0.002217805172987552
$ws.Range("C36").Value = 0.007315981828670959
$ws.Range("D36").Value = 0.002452849281681767
$ws.Range("C37").Value = 0.0078885927962675
$ws.Range("D37").Value = 0.002703076696549811
$ws.Range("C38").Value = 0.008488512532464746
$ws.Range("D38").Value = 0.002968818256483796
$ws.Range("C39").Value = 0.009117915647600672
$ws.Range("D39").Value = 0.003250325080223701
$ws.Range("C40").Value = 0.009775566679519307
$ws.Range("D40").Value = 0.003547496145192132
$ws.Range("C41").Value = 0.01046174078788672
$ws.Range("D41").Value = 0.003860251135713529
$ws.Range("C42").Value = 0.01117585559508204
$ws.Range("D42").Value = 0.004188726913248905
$ws.Range("C43").Value = 0.01192132392023616
$ws.Range("D43").Value = 0.004534531481107271
$ws.Range("C44").Value = 0.01269854245234837
$ws.Range("D44").Value = 0.004898019053775293
$ws.Range("C45").Value = 0.01350844549611614
$ws.Range("D45").Value = 0.00527870927697528
$ws.Range("C46").Value = 0.014352092953841
$ws.Range("D46").Value = 0.00567925265883003
$ws.Range("C47").Value = 0.01523266151243704
$ws.Range("D47").Value = 0.006098057463777794
$ws.Range("C48").Value = 0.01614964635304865
$ws.Range("D48").Value = 0.006537030950276152
$ws.Range("C49").Value = 0.01710403762537641
$ws.Range("D49").Value = 0.006998998747128493
$ws.Range("C50").Value = 0.01809596405071672
$ws.Range("D50").Value = 0.007483677327515616
$ws.Range("C51").Value = 0.01912612660035589
$ws.Range("D51").Value = 0.007990666079910532
$ws.Range("C52").Value = 0.02019344679425151
$ws.Range("D52").Value = 0.008519422283512317
$ws.Range("C53").Value = 0.02130386235023082
$ws.Range("D53").Value = 0.009070982358277729
$ws.Range("C54").Value = 0.02246009240890871
$ws.Range("D54").Value = 0.009647720499109285
$ws.Range("C55").Value = 0.02366191496656731
$ws.Range("D55").Value = 0.01025070962051034
$ws.Range("C56").Value = 0.02491145968664366
$ws.Range("D56").Value = 0.01088081613215369
$ws.Range("C57").Value = 0.02621089514218796
$ws.Range("D57").Value = 0.01153878917353518
$ws.Range("C58").Value = 0.02756270910947833
$ws.Range("D58").Value = 0.01222621898230268
$ws.Range("C59").Value = 0.02897082583611871
$ws.Range("D59").Value = 0.01294195311680361
$ws.Range("C60").Value = 0.03043563762988971
$ws.Range("D60").Value = 0.01368522223732832
$ws.Range("C61").Value = 0.03195584506948419
$ws.Range("D61").Value = 0.01445919199768124
$ws.Range("C62").Value = 0.03353463008261574
$ws.Range("D62").Value = 0.01526518928522395
$ws.Range("C63").Value = 0.03517350216937644
$ws.Range("D63").Value = 0.01610568252518509
$ws.Range("C64").Value = 0.03688458123279914
$ws.Range("D64").Value = 0.0169791448371517
$ws.Range("C65").Value = 0.03867490910589202
$ws.Range("D65").Value = 0.01788851024295812
$ws.Range("C66").Value = 0.04054354385268845
$ws.Range("D66").Value = 0.01883536120312562
$ws.Range("C67").Value = 0.04250457018232223
$ws.Range("D67").Value = 0.01981957949496243
$ws.Range("C68").Value = 0.04456492103503033
$ws.Range("D68").Value = 0.020843863502726
$ws.Range("C69").Value = 0.0467278280438068
$ws.Range("D69").Value = 0.02191041957886909
$ws.Range("C70").Value = 0.04900743513945455
$ws.Range("D70").Value = 0.02302179042945563
$ws.Range("C71").Value = 0.05141787819337239
$ws.Range("D71").Value = 0.02418256334479382
$ws.Range("C72").Value = 0.05398091831767622
$ws.Range("D72").Value = 0.02539032165481049
$ws.Range("C73").Value = 0.05671953637977685
$ws.Range("D73").Value = 0.02664892942348986
$ws.Range("C74").Value = 0.05965117860921963
$ws.Range("D74").Value = 0.02796490881600504
$ws.Range("C75").Value = 0.06279441370588512
$ws.Range("D75").Value = 0.02934538623441432
$ws.Range("C76").Value = 0.0661837065941792
$ws.Range("D76").Value = 0.03079981660540627
$ws.Range("C77").Value = 0.06990208124469079
$ws.Range("D77").Value = 0.03232659487315364
$ws.Range("C78").Value = 0.07400343836646275
$ws.Range("D78").Value = 0.03394091484267808
$ws.Range("C79").Value = 0.07855896095546797
$ws.Range("D79").Value = 0.03564050952743403
$ws.Range("C80").Value = 0.08362345145617162
$ws.Range("D80").Value = 0.037444889464284
$ws.Range("C81").Value = 0.08925581334083788
$ws.Range("D81").Value = 0.03936855756950733
$ws.Range("C82").Value = 0.09548072573510334
$ws.Range("D82").Value = 0.04143620376145576
$ws.Range("C83").Value = 0.1023557961971335
$ws.Range("D83").Value = 0.04367057870975791
$ws.Range("C84").Value = 0.1099821286977521
$ws.Range("D84").Value = 0.0461224126985918
$ws.Range("C85").Value = 0.1184396291548474
$ws.Range("D85").Value = 0.04883954865913382
$ws.Range("C86").Value = 0.12790343816437
$ws.Range("D86").Value = 0.05194420497538475
$ws.Range("C87").Value = 0.1385371133944896
$ws.Range("D87").Value = 0.05555706999916801
$ws.Range("C88").Value = 0.1507233118955087
$ws.Range("D88").Value = 0.06011473886522573
$ws.Range("C89").Value = 0.1649200401693039
$ws.Range("D89").Value = 0.0676966838447788
$ws.Range("C90").Value = 0.1822302172879039
$ws.Range("D90").Value = 0.08333357312648947
$ws.Range("C91").Value = 0.2047957061264629
$ws.Range("D91").Value = 0.1071827620336817
$ws.Range("C92").Value = 0.2341845279079098
$ws.Range("D92").Value = 0.1392612400619395
$ws.Range("C93").Value = 0.271341383841199
$ws.Range("D93").Value = 0.1800641854559589
$ws.Range("C94").Value = 0.3158288464589407
$ws.Range("D94").Value = 0.2290100809710532
$ws.Range("C95").Value = 0.3672233655101436
$ws.Range("D95").Value = 0.2859495174893664
$ws.Range("C96").Value = 0.4257459704027849
$ws.Range("D96").Value = 0.3510517528023256
$ws.Range("C97").Value = 0.4921630365870894
$ws.Range("D97").Value = 0.4250329717688658
$ws.Range("C98").Value = 0.5659931523888523
$ws.Range("D98").Value = 0.5077496658443149
$ws.Range("C99").Value = 0.6483115577321803
$ws.Range("D99").Value = 0.6003524429479112
$ws.Range("C100").Value = 0.7410554761467684
$ws.Range("D100").Value = 0.7050360616910133
$ws.Range("C101").Value = 0.8489297142699497
$ws.Range("D101").Value = 0.827393628201106
